# Update omics and assembly tags on the isa_template sheet.
# The "Tags" row (row 15) previously tagged this template with "genome" (EFO:0004420),
# it is now tagged "Genomics" (NCIT_C84343) instead, and the matching "Tags Term Source REF"
# (row 17) changes from EFO to NCIT to match the new term source.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

$ws.Range("B15").Value = "Genomics"
$ws.Range("B16").Value = "http://purl.obolibrary.org/obo/NCIT_C84343"
$ws.Range("B17").Value = "NCIT"
